$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new data row (row 2) of the Campos Register datapool.
$ws.Range("A2").Value = "Fabian"
$ws.Range("B2").Value = "Alfonso"
$ws.Range("C2").Value = 311
$ws.Range("D2").Value = "elkin3001"
$ws.Range("E2").Value = "Galan M#14"
$ws.Range("F2").Value = "Calarca"
$ws.Range("G2").Value = "Quindio"
$ws.Range("H2").Value = 57
$ws.Range("I2").Value = "Colombia"
$ws.Range("J2").Value = "elkin3001"
$ws.Range("K2").Value = "elkin3001"
$ws.Range("L2").Value = "elkin3001"

# Move/leave the active selection where Excel landed after data entry.
$ws.Range("L3").Select()
